$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.287.37'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.26%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.881.74'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.92%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.54'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.68'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.11%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.881.16'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.89'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -5.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.145'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.431'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000232'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.08%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.86%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.349.45'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.283.83'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.19%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.886.57'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '431.46'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.05'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.18%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.82'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.73%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.25'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.86'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.01'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -9.64%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000104'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.99'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.26%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.06'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -7.94%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.46'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.958'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.38'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.04%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.63%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -8.09%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.91%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.71%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '38.97'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.706.31'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '132.92'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0336'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '342.01'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.18%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.52'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.84%  '
